$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Block 1 (fcn00): row 4 - new NNLAB simulated-corneal-reflection estimates ---
$ws.Range("F4").Value = 0.56208599999999997
$ws.Range("G4").Value = 2.3933610000000001
$ws.Range("H4").Value = -0.56828199999999995
$ws.Range("I4").Value = "NNLAB: Estimated 41/48"

# --- Block 2 (fcn01): row 25 ---
$ws.Range("F25").Value = 0.54704799999999998
$ws.Range("G25").Value = 2.2653979999999998
$ws.Range("H25").Value = -0.62666500000000003
$ws.Range("I25").Value = "NNLAB: estimated 39/48"

# --- Block 3 (fcn02): row 45 ---
$ws.Range("F45").Value = 0.481323
$ws.Range("G45").Value = 2.0770569999999999
$ws.Range("H45").Value = -0.62484499999999998
$ws.Range("I45").Value = "NNLAB: Estimated 39 / 48"

# --- Block 3 hidden/centering rows: row 57 ---
$ws.Range("G57").Value = 2.754143
$ws.Range("H57").Value = -0.65172699999999995

# --- Sheet2 mirrored summary table: row 66 ---
$ws.Range("F66").Value = 0.481323
$ws.Range("G66").Value = 2.754143
$ws.Range("H66").Value = -0.65172699999999995
$ws.Range("I66").Value = "NNLAB: Estimated 36/48"

# --- View adjustments: scroll right a bit and move the selection ---
$ws.Columns.Item(5).ColumnWidth = 14

$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$ws.Range("I5").Select()
